$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, DateSerial, Interval, Prediction, Lookup
$data = @(
@(2,45890,17,0,'21.08.202517'),
@(3,45890,18,0,'21.08.202518'),
@(4,45890,19,0.417,'21.08.202519'),
@(5,45890,20,0.053,'21.08.202520'),
@(6,45890,21,0,'21.08.202521'),
@(7,45890,22,0,'21.08.202522'),
@(8,45890,23,0,'21.08.202523'),
@(9,45890,24,0,'21.08.202524'),
@(10,45891,1,0,'22.08.20251'),
@(11,45891,2,0,'22.08.20252'),
@(12,45891,3,0,'22.08.20253'),
@(13,45891,4,0,'22.08.20254'),
@(14,45891,5,0,'22.08.20255'),
@(15,45891,6,0.0,'22.08.20256'),
@(16,45891,7,0.0,'22.08.20257'),
@(17,45891,8,0.153,'22.08.20258'),
@(18,45891,9,0.689,'22.08.20259'),
@(19,45891,10,1.214,'22.08.202510'),
@(20,45891,11,1.95,'22.08.202511'),
@(21,45891,12,2.237,'22.08.202512'),
@(22,45891,13,1.997,'22.08.202513'),
@(23,45891,14,1.977,'22.08.202514'),
@(24,45891,15,1.502,'22.08.202515'),
@(25,45891,16,1.036,'22.08.202516'),
@(26,45891,17,0.932,'22.08.202517'),
@(27,45891,18,0.68,'22.08.202518'),
@(28,45891,19,0.255,'22.08.202519'),
@(29,45891,20,0.044,'22.08.202520'),
@(30,45891,21,0,'22.08.202521'),
@(31,45891,22,0,'22.08.202522'),
@(32,45891,23,0,'22.08.202523'),
@(33,45891,24,0,'22.08.202524'),
@(34,45892,1,0,'23.08.20251'),
@(35,45892,2,0,'23.08.20252'),
@(36,45892,3,0,'23.08.20253'),
@(37,45892,4,0,'23.08.20254'),
@(38,45892,5,0,'23.08.20255'),
@(39,45892,6,0.0,'23.08.20256'),
@(40,45892,7,0.0,'23.08.20257'),
@(41,45892,8,0.267,'23.08.20258'),
@(42,45892,9,1.061,'23.08.20259'),
@(43,45892,10,1.74,'23.08.202510'),
@(44,45892,11,2.53,'23.08.202511'),
@(45,45892,12,3.04,'23.08.202512'),
@(46,45892,13,3.251,'23.08.202513'),
@(47,45892,14,3.197,'23.08.202514'),
@(48,45892,15,3.044,'23.08.202515'),
@(49,45892,16,2.691,'23.08.202516'),
@(50,45892,17,1.993,'23.08.202517'),
@(51,45892,18,1.23,'23.08.202518'),
@(52,45892,19,0.529,'23.08.202519'),
@(53,45892,20,0.117,'23.08.202520'),
@(54,45892,21,0,'23.08.202521'),
@(55,45892,22,0,'23.08.202522'),
@(56,45892,23,0,'23.08.202523'),
@(57,45892,24,0,'23.08.202524'),
@(58,45893,1,0,'24.08.20251'),
@(59,45893,2,0,'24.08.20252'),
@(60,45893,3,0,'24.08.20253'),
@(61,45893,4,0,'24.08.20254'),
@(62,45893,5,0,'24.08.20255'),
@(63,45893,6,0.0,'24.08.20256'),
@(64,45893,7,0.0,'24.08.20257'),
@(65,45893,8,0.272,'24.08.20258'),
@(66,45893,9,1.02,'24.08.20259'),
@(67,45893,10,1.696,'24.08.202510'),
@(68,45893,11,2.293,'24.08.202511'),
@(69,45893,12,2.784,'24.08.202512'),
@(70,45893,13,2.868,'24.08.202513'),
@(71,45893,14,2.901,'24.08.202514'),
@(72,45893,15,2.859,'24.08.202515'),
@(73,45893,16,2.369,'24.08.202516'),
@(74,45893,17,2.151,'24.08.202517'),
@(75,45893,18,1.416,'24.08.202518'),
@(76,45893,19,0.69,'24.08.202519'),
@(77,45893,20,0.15,'24.08.202520'),
@(78,45893,21,0,'24.08.202521'),
@(79,45893,22,0,'24.08.202522'),
@(80,45893,23,0,'24.08.202523'),
@(81,45893,24,0,'24.08.202524'),
@(82,45894,1,0,'25.08.20251'),
@(83,45894,2,0,'25.08.20252'),
@(84,45894,3,0,'25.08.20253'),
@(85,45894,4,0,'25.08.20254'),
@(86,45894,5,0,'25.08.20255'),
@(87,45894,6,0.0,'25.08.20256'),
@(88,45894,7,0.0,'25.08.20257'),
@(89,45894,8,0.3,'25.08.20258'),
@(90,45894,9,1.296,'25.08.20259'),
@(91,45894,10,2.359,'25.08.202510'),
@(92,45894,11,3.016,'25.08.202511'),
@(93,45894,12,3.166,'25.08.202512'),
@(94,45894,13,3.357,'25.08.202513'),
@(95,45894,14,3.388,'25.08.202514'),
@(96,45894,15,3.315,'25.08.202515'),
@(97,45894,16,3.034,'25.08.202516'),
@(98,45894,17,2.483,'25.08.202517'),
@(99,45894,18,1.709,'25.08.202518'),
@(100,45894,19,0.748,'25.08.202519'),
@(101,45894,20,0.142,'25.08.202520'),
@(102,45894,21,0,'25.08.202521'),
@(103,45894,22,0,'25.08.202522'),
@(104,45894,23,0,'25.08.202523'),
@(105,45894,24,0,'25.08.202524'),
@(106,45895,1,0,'26.08.20251'),
@(107,45895,2,0,'26.08.20252'),
@(108,45895,3,0,'26.08.20253'),
@(109,45895,4,0,'26.08.20254'),
@(110,45895,5,0,'26.08.20255'),
@(111,45895,6,0.0,'26.08.20256'),
@(112,45895,7,0.0,'26.08.20257'),
@(113,45895,8,0.289,'26.08.20258'),
@(114,45895,9,1.208,'26.08.20259'),
@(115,45895,10,2.049,'26.08.202510'),
@(116,45895,11,2.508,'26.08.202511'),
@(117,45895,12,2.989,'26.08.202512'),
@(118,45895,13,3.177,'26.08.202513'),
@(119,45895,14,3.16,'26.08.202514'),
@(120,45895,15,3.051,'26.08.202515'),
@(121,45895,16,2.794,'26.08.202516'),
@(122,45895,17,2.322,'26.08.202517'),
@(123,45895,18,1.649,'26.08.202518'),
@(124,45895,19,0.827,'26.08.202519'),
@(125,45895,20,0.158,'26.08.202520'),
@(126,45895,21,0,'26.08.202521'),
@(127,45895,22,0,'26.08.202522'),
@(128,45895,23,0,'26.08.202523'),
@(129,45895,24,0,'26.08.202524'),
@(130,45896,1,0,'27.08.20251'),
@(131,45896,2,0,'27.08.20252'),
@(132,45896,3,0,'27.08.20253'),
@(133,45896,4,0,'27.08.20254'),
@(134,45896,5,0,'27.08.20255'),
@(135,45896,6,0.0,'27.08.20256'),
@(136,45896,7,0.0,'27.08.20257'),
@(137,45896,8,0.204,'27.08.20258'),
@(138,45896,9,0.817,'27.08.20259'),
@(139,45896,10,1.451,'27.08.202510'),
@(140,45896,11,2.001,'27.08.202511'),
@(141,45896,12,2.359,'27.08.202512'),
@(142,45896,13,2.446,'27.08.202513'),
@(143,45896,14,2.454,'27.08.202514'),
@(144,45896,15,2.312,'27.08.202515'),
@(145,45896,16,2.143,'27.08.202516'),
@(146,45896,17,1.718,'27.08.202517'),
@(147,45896,18,1.117,'27.08.202518'),
@(148,45896,19,0.516,'27.08.202519'),
@(149,45896,20,0.108,'27.08.202520'),
@(150,45896,21,0,'27.08.202521'),
@(151,45896,22,0,'27.08.202522'),
@(152,45896,23,0,'27.08.202523'),
@(153,45896,24,0,'27.08.202524'),
@(154,45897,1,0,'28.08.20251'),
@(155,45897,2,0,'28.08.20252'),
@(156,45897,3,0,'28.08.20253'),
@(157,45897,4,0,'28.08.20254'),
@(158,45897,5,0,'28.08.20255'),
@(159,45897,6,0.0,'28.08.20256'),
@(160,45897,7,0.0,'28.08.20257'),
@(161,45897,8,0.288,'28.08.20258'),
@(162,45897,9,1.108,'28.08.20259'),
@(163,45897,10,1.992,'28.08.202510'),
@(164,45897,11,2.51,'28.08.202511'),
@(165,45897,12,2.816,'28.08.202512'),
@(166,45897,13,3.079,'28.08.202513'),
@(167,45897,14,3.16,'28.08.202514'),
@(168,45897,15,3.032,'28.08.202515'),
@(169,45897,16,2.847,'28.08.202516'),
@(170,45897,17,2.302,'28.08.202517')
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}

Write-Host "Updated $($data.Count) rows"
